# Auto-generated edit script: apply scheduled-runner price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 11173.071
$ws.Cells.Item(86, 9).Value = 15482.875
$ws.Cells.Item(86, 10).Value = 5426.6665
$ws.Cells.Item(86, 11).Value = 15482.875
$ws.Cells.Item(86, 12).Value = 5426.6665
$ws.Cells.Item(86, 13).Value = -14359.875
$ws.Cells.Item(86, 14).Value = -7672.6665
$ws.Cells.Item(89, 8).Value = 11173.071
$ws.Cells.Item(89, 9).Value = 15482.875
$ws.Cells.Item(89, 10).Value = 5426.6665
$ws.Cells.Item(89, 11).Value = 77414.375
$ws.Cells.Item(89, 12).Value = 27133.3325
$ws.Cells.Item(89, 13).Value = -71798.375
$ws.Cells.Item(89, 14).Value = -38365.3325
$ws.Cells.Item(96, 8).Value = 374.2414
$ws.Cells.Item(96, 9).Value = 291.65
$ws.Cells.Item(96, 10).Value = 557.7778
$ws.Cells.Item(96, 11).Value = 874.9499999999999
$ws.Cells.Item(96, 12).Value = 1673.3334
$ws.Cells.Item(96, 13).Value = 498.0500000000001
$ws.Cells.Item(96, 14).Value = -4419.3334
$ws.Cells.Item(100, 8).Value = 31251756
$ws.Cells.Item(100, 9).Value = 1810.5
$ws.Cells.Item(100, 10).Value = 83335000
$ws.Cells.Item(100, 11).Value = 1810.5
$ws.Cells.Item(100, 12).Value = 83335000
$ws.Cells.Item(100, 13).Value = -1269.5
$ws.Cells.Item(100, 14).Value = -83336082
$ws.Cells.Item(101, 8).Value = 1103.3478
$ws.Cells.Item(101, 9).Value = 965.125
$ws.Cells.Item(101, 10).Value = 1419.2858
$ws.Cells.Item(101, 11).Value = 2895.375
$ws.Cells.Item(101, 12).Value = 4257.857400000001
$ws.Cells.Item(101, 13).Value = -1273.375
$ws.Cells.Item(101, 14).Value = -7501.857400000001
$ws.Cells.Item(137, 8).Value = 4681.8184
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1561.5428
$ws.Cells.Item(74, 9).Value = 716.7879
$ws.Cells.Item(74, 10).Value = 15500
$ws.Cells.Item(74, 11).Value = 716.7879
$ws.Cells.Item(74, 12).Value = 15500
$ws.Cells.Item(74, 13).Value = 157.2121
$ws.Cells.Item(74, 14).Value = -17248
$ws.Cells.Item(77, 8).Value = 1561.5428
$ws.Cells.Item(77, 9).Value = 716.7879
$ws.Cells.Item(77, 10).Value = 15500
$ws.Cells.Item(77, 11).Value = 3583.9395
$ws.Cells.Item(77, 12).Value = 77500
$ws.Cells.Item(77, 13).Value = 784.0604999999996
$ws.Cells.Item(77, 14).Value = -86236
$ws.Cells.Item(121, 8).Value = 32121
$ws.Cells.Item(121, 10).Value = 32121
$ws.Cells.Item(121, 12).Value = 32121
$ws.Cells.Item(121, 14).Value = -35615
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 22422.125
$ws.Cells.Item(82, 9).Value = 10325.167
$ws.Cells.Item(82, 10).Value = 29680.3
$ws.Cells.Item(82, 11).Value = 10325.167
$ws.Cells.Item(82, 12).Value = 29680.3
$ws.Cells.Item(82, 13).Value = -9942.166999999999
$ws.Cells.Item(82, 14).Value = -30446.3
$ws.Cells.Item(85, 8).Value = 22422.125
$ws.Cells.Item(85, 9).Value = 10325.167
$ws.Cells.Item(85, 10).Value = 29680.3
$ws.Cells.Item(85, 11).Value = 10325.167
$ws.Cells.Item(85, 12).Value = 29680.3
$ws.Cells.Item(85, 13).Value = -8999.166999999999
$ws.Cells.Item(85, 14).Value = -32332.3
$ws.Cells.Item(86, 8).Value = 1665.119
$ws.Cells.Item(86, 9).Value = 1716.5862
$ws.Cells.Item(86, 10).Value = 1550.3077
$ws.Cells.Item(86, 11).Value = 1716.5862
$ws.Cells.Item(86, 12).Value = 1550.3077
$ws.Cells.Item(86, 13).Value = -593.5862
$ws.Cells.Item(86, 14).Value = -3796.3077
$ws.Cells.Item(89, 8).Value = 1665.119
$ws.Cells.Item(89, 9).Value = 1716.5862
$ws.Cells.Item(89, 10).Value = 1550.3077
$ws.Cells.Item(89, 11).Value = 8582.931
$ws.Cells.Item(89, 12).Value = 7751.538500000001
$ws.Cells.Item(89, 13).Value = -2966.931
$ws.Cells.Item(89, 14).Value = -18983.5385
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5113
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 5113
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).ClearContents()
$ws.Cells.Item(31, 13).Value = 5113
$ws.Cells.Item(31, 14).Value = -5703
$ws.Cells.Item(34, 8).Value = 5113
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 5113
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).ClearContents()
$ws.Cells.Item(34, 13).Value = 5113
$ws.Cells.Item(34, 14).Value = -5517
$ws.Cells.Item(41, 8).Value = 15282.5
$ws.Cells.Item(41, 9).Value = 3950
$ws.Cells.Item(41, 11).Value = 3950
$ws.Cells.Item(41, 13).Value = -3522
$ws.Cells.Item(50, 8).Value = 9099.833000000001
$ws.Cells.Item(50, 10).Value = 9099.833000000001
$ws.Cells.Item(50, 12).Value = 9099.833000000001
$ws.Cells.Item(50, 14).Value = -10349.833
$ws.Cells.Item(51, 8).Value = 9132.75
$ws.Cells.Item(51, 10).Value = 9132.75
$ws.Cells.Item(51, 12).Value = 9132.75
$ws.Cells.Item(51, 14).Value = -10604.75
$ws.Cells.Item(60, 8).Value = 23218.076
$ws.Cells.Item(60, 9).Value = 93
$ws.Cells.Item(60, 10).Value = 25145.166
$ws.Cells.Item(60, 11).Value = 93
$ws.Cells.Item(60, 12).Value = 25145.166
$ws.Cells.Item(60, 13).Value = 418
$ws.Cells.Item(60, 14).Value = -26167.166
$ws.Cells.Item(61, 8).Value = 9132.75
$ws.Cells.Item(61, 10).Value = 9132.75
$ws.Cells.Item(61, 12).Value = 9132.75
$ws.Cells.Item(61, 14).Value = -9828.75
$ws.Cells.Item(68, 8).Value = 17500
$ws.Cells.Item(68, 10).Value = 17500
$ws.Cells.Item(68, 12).Value = 17500
$ws.Cells.Item(68, 14).Value = -18998
$ws.Cells.Item(71, 8).Value = 17500
$ws.Cells.Item(71, 10).Value = 17500
$ws.Cells.Item(71, 12).Value = 52500
$ws.Cells.Item(71, 14).Value = -59988
$ws.Cells.Item(99, 8).Value = 1576.3334
$ws.Cells.Item(99, 9).Value = 1491.2727
$ws.Cells.Item(99, 10).Value = 1710
$ws.Cells.Item(99, 11).Value = 1491.2727
$ws.Cells.Item(99, 12).Value = 1710
$ws.Cells.Item(99, 13).Value = 6.727300000000014
$ws.Cells.Item(99, 14).Value = -4706
$ws.Cells.Item(109, 8).Value = 11080
$ws.Cells.Item(109, 10).Value = 11080
$ws.Cells.Item(109, 12).Value = 11080
$ws.Cells.Item(109, 14).Value = -13160
$ws.Cells.Item(126, 8).Value = 1576.3334
$ws.Cells.Item(126, 9).Value = 1491.2727
$ws.Cells.Item(126, 10).Value = 1710
$ws.Cells.Item(126, 11).Value = 4473.8181
$ws.Cells.Item(126, 12).Value = 5130
$ws.Cells.Item(126, 13).Value = -2003.8181
$ws.Cells.Item(126, 14).Value = -10070
$ws.Cells.Item(132, 8).Value = 1663.3793
$ws.Cells.Item(132, 9).Value = 1827.8096
$ws.Cells.Item(132, 10).Value = 1231.75
$ws.Cells.Item(132, 11).Value = 5483.4288
$ws.Cells.Item(132, 12).Value = 3695.25
$ws.Cells.Item(132, 13).Value = -2953.4288
$ws.Cells.Item(132, 14).Value = -8755.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 75.875
$ws.Cells.Item(8, 9).Value = 75.875
$ws.Cells.Item(8, 11).Value = 227.625
$ws.Cells.Item(8, 13).Value = -88.625
$ws.Cells.Item(43, 8).Value = 7500
$ws.Cells.Item(43, 10).Value = 7500
$ws.Cells.Item(43, 12).Value = 22500
$ws.Cells.Item(43, 14).Value = -22728
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(51, 8).Value = 44980
$ws.Cells.Item(51, 10).Value = 44980
$ws.Cells.Item(51, 12).Value = 44980
$ws.Cells.Item(51, 14).Value = -45998
$ws.Cells.Item(57, 8).Value = 12957.6
$ws.Cells.Item(57, 10).Value = 14933.25
$ws.Cells.Item(57, 12).Value = 14933.25
$ws.Cells.Item(57, 14).Value = -16573.25
$ws.Cells.Item(122, 8).Value = 1839.9395
$ws.Cells.Item(122, 9).Value = 1797.4375
$ws.Cells.Item(122, 10).Value = 3200
$ws.Cells.Item(122, 11).Value = 5392.3125
$ws.Cells.Item(122, 12).Value = 9600
$ws.Cells.Item(122, 13).Value = -2942.3125
$ws.Cells.Item(122, 14).Value = -14500
$ws.Cells.Item(123, 8).Value = 39509.43
$ws.Cells.Item(123, 10).Value = 39509.43
$ws.Cells.Item(123, 12).Value = 39509.43
$ws.Cells.Item(123, 14).Value = -44409.43
$ws.Cells.Item(132, 8).Value = 1459.5082
$ws.Cells.Item(132, 9).Value = 1150.3191
$ws.Cells.Item(132, 11).Value = 3450.9573
$ws.Cells.Item(132, 13).Value = -920.9573
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6697.6875
$ws.Cells.Item(132, 9).Value = 7124.1816
$ws.Cells.Item(132, 10).Value = 5759.4
$ws.Cells.Item(132, 11).Value = 21372.5448
$ws.Cells.Item(132, 12).Value = 17278.2
$ws.Cells.Item(132, 13).Value = -18842.5448
$ws.Cells.Item(132, 14).Value = -22338.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 26200
$ws.Cells.Item(16, 10).Value = 26200
$ws.Cells.Item(16, 12).Value = 26200
$ws.Cells.Item(16, 14).Value = -26784
$ws.Cells.Item(100, 8).Value = 1112109.5
$ws.Cells.Item(100, 9).Value = 1023
$ws.Cells.Item(100, 10).Value = 2000978.6
$ws.Cells.Item(100, 11).Value = 2046
$ws.Cells.Item(100, 12).Value = 4001957.2
$ws.Cells.Item(100, 13).Value = -1505
$ws.Cells.Item(100, 14).Value = -4003039.2
$ws.Cells.Item(109, 8).Value = 27888.5
$ws.Cells.Item(109, 10).Value = 27888.5
$ws.Cells.Item(109, 12).Value = 27888.5
$ws.Cells.Item(109, 14).Value = -30662.5
$ws.Cells.Item(132, 8).Value = 24797.227
$ws.Cells.Item(132, 9).Value = 2855.7727
$ws.Cells.Item(132, 10).Value = 46738.684
$ws.Cells.Item(132, 11).Value = 8567.3181
$ws.Cells.Item(132, 12).Value = 140216.052
$ws.Cells.Item(132, 13).Value = -6037.3181
$ws.Cells.Item(132, 14).Value = -145276.052
